$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H43").Value = 8481
$ws.Range("J43").Value = 8101
$ws.Range("L43").Value = 8101
$ws.Range("N43").Value = -8239
$ws.Range("H98").Value = 1903.6129
$ws.Range("J98").Value = 1461.75
$ws.Range("L98").Value = 1461.75
$ws.Range("N98").Value = -4457.75
$ws.Range("H105").Value = 85000
$ws.Range("J105").Value = 85000
$ws.Range("L105").Value = 85000
$ws.Range("N105").Value = -91988
$ws.Range("H106").Value = 2992.5557
$ws.Range("I106").Value = 2979.5
$ws.Range("J106").Value = 3097
$ws.Range("K106").Value = 2979.5
$ws.Range("L106").Value = 3097
$ws.Range("M106").Value = -2348.5
$ws.Range("N106").Value = -4359
$ws.Range("H107").Value = 1955.619
$ws.Range("I107").Value = 1734.9412
$ws.Range("J107").Value = 2893.5
$ws.Range("K107").Value = 1734.9412
$ws.Range("L107").Value = 2893.5
$ws.Range("M107").Value = 185.0588
$ws.Range("N107").Value = -6733.5
$ws.Range("H116").Value = 16579.75
$ws.Range("I116").Value = 25924.445
$ws.Range("J116").Value = 8934.091
$ws.Range("K116").Value = 25924.445
$ws.Range("L116").Value = 8934.091
$ws.Range("M116").Value = -22482.445
$ws.Range("N116").Value = -15818.091
$ws.Range("H122").Value = 1903.6129
$ws.Range("J122").Value = 1461.75
$ws.Range("L122").Value = 4385.25
$ws.Range("N122").Value = -9285.25
$ws.Range("H137").Value = 8763.578
$ws.Range("I137").Value = 12146.586
$ws.Range("K137").Value = 36439.758
$ws.Range("M137").Value = -33889.758
$ws.Range("H138").Value = 29446.865
$ws.Range("I138").Value = 1908.6086
$ws.Range("J138").Value = 74688.28999999999
$ws.Range("K138").Value = 5725.825800000001
$ws.Range("L138").Value = 224064.87
$ws.Range("M138").Value = -585.8258000000005
$ws.Range("N138").Value = -234344.87

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5219.5713
$ws.Range("I28").Value = 5219.5713
$ws.Range("K28").Value = 5219.5713
$ws.Range("M28").Value = -5027.5713
$ws.Range("H45").Value = 2197.3215
$ws.Range("I45").Value = 1383
$ws.Range("K45").Value = 1383
$ws.Range("M45").Value = -1006
$ws.Range("H74").Value = 266845.6
$ws.Range("J74").Value = 15506.875
$ws.Range("L74").Value = 15506.875
$ws.Range("N74").Value = -17254.875
$ws.Range("H77").Value = 266845.6
$ws.Range("J77").Value = 15506.875
$ws.Range("L77").Value = 77534.375
$ws.Range("N77").Value = -86270.375
$ws.Range("H99").Value = 5219.5713
$ws.Range("I99").Value = 5219.5713
$ws.Range("K99").Value = 5219.5713
$ws.Range("M99").Value = -2224.5713
$ws.Range("H122").Value = 1775.2188
$ws.Range("I122").Value = 1760.5333
$ws.Range("K122").Value = 5281.5999
$ws.Range("M122").Value = -2831.5999
$ws.Range("H137").Value = 118674.375
$ws.Range("J137").Value = 118674.375
$ws.Range("L137").Value = 118674.375
$ws.Range("N137").Value = -128874.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15221.087
$ws.Range("I20").Value = 19923.059
$ws.Range("J20").Value = 1898.8334
$ws.Range("K20").Value = 19923.059
$ws.Range("L20").Value = 1898.8334
$ws.Range("M20").Value = -19676.059
$ws.Range("N20").Value = -2392.8334
$ws.Range("H86").Value = 1742.92
$ws.Range("I86").Value = 1760.8235
$ws.Range("J86").Value = 1704.875
$ws.Range("K86").Value = 1760.8235
$ws.Range("L86").Value = 1704.875
$ws.Range("M86").Value = -637.8235
$ws.Range("N86").Value = -3950.875
$ws.Range("H89").Value = 1742.92
$ws.Range("I89").Value = 1760.8235
$ws.Range("J89").Value = 1704.875
$ws.Range("K89").Value = 8804.1175
$ws.Range("L89").Value = 8524.375
$ws.Range("M89").Value = -3188.1175
$ws.Range("N89").Value = -19756.375
$ws.Range("H99").Value = 1929.5883
$ws.Range("J99").Value = 4007
$ws.Range("L99").Value = 4007
$ws.Range("N99").Value = -7003
$ws.Range("H105").Value = 3753.2727
$ws.Range("I105").Value = 1854.75
$ws.Range("K105").Value = 1854.75
$ws.Range("M105").Value = -107.75
$ws.Range("H134").Value = 1873.1915
$ws.Range("I134").Value = 1464.575
$ws.Range("K134").Value = 4393.725
$ws.Range("M134").Value = -1858.725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -112488
$ws.Range("H99").Value = 5620.1577
$ws.Range("I99").Value = 3884.9333
$ws.Range("J99").Value = 12127.25
$ws.Range("K99").Value = 3884.9333
$ws.Range("L99").Value = 12127.25
$ws.Range("M99").Value = -2386.9333
$ws.Range("N99").Value = -15123.25
$ws.Range("H126").Value = 5620.1577
$ws.Range("I126").Value = 3884.9333
$ws.Range("J126").Value = 12127.25
$ws.Range("K126").Value = 11654.7999
$ws.Range("L126").Value = 36381.75
$ws.Range("M126").Value = -9184.7999
$ws.Range("N126").Value = -41321.75
$ws.Range("H140").Value = 119899
$ws.Range("J140").Value = 119899
$ws.Range("L140").Value = 119899
$ws.Range("N140").Value = -130259

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2184.8
$ws.Range("J6").Value = 10
$ws.Range("L6").Value = 30
$ws.Range("N6").Value = -256
$ws.Range("H58").Value = 3424.6155
$ws.Range("I58").Value = 1125
$ws.Range("J58").Value = 4446.6665
$ws.Range("K58").Value = 3375
$ws.Range("L58").Value = 13339.9995
$ws.Range("M58").Value = -3247
$ws.Range("N58").Value = -13595.9995
$ws.Range("H121").Value = 43764.32
$ws.Range("I121").Value = 89429.664
$ws.Range("K121").Value = 268288.992
$ws.Range("M121").Value = -266978.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11607.333
$ws.Range("I80").Value = 2799.5715
$ws.Range("J80").Value = 17212.273
$ws.Range("K80").Value = 2799.5715
$ws.Range("L80").Value = 17212.273
$ws.Range("M80").Value = -1801.5715
$ws.Range("N80").Value = -19208.273
$ws.Range("H83").Value = 11607.333
$ws.Range("I83").Value = 2799.5715
$ws.Range("J83").Value = 17212.273
$ws.Range("K83").Value = 13997.8575
$ws.Range("L83").Value = 86061.36500000001
$ws.Range("M83").Value = -9005.8575
$ws.Range("N83").Value = -96045.36500000001
$ws.Range("H99").Value = 7585.4287
$ws.Range("I99").Value = 7585.4287
$ws.Range("K99").Value = 7585.4287
$ws.Range("M99").Value = -5339.4287
$ws.Range("H102").Value = 5097.7666
$ws.Range("I102").Value = 5219.0356
$ws.Range("K102").Value = 5219.0356
$ws.Range("M102").Value = -3597.0356
$ws.Range("H136").Value = 8998.75
$ws.Range("J136").Value = 8998.75
$ws.Range("L136").Value = 26996.25
$ws.Range("N136").Value = -32096.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1300.8182
$ws.Range("I132").Value = 1000.9
$ws.Range("K132").Value = 3002.7
$ws.Range("M132").Value = -472.6999999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2266.4736
$ws.Range("I132").Value = 1006.9
$ws.Range("K132").Value = 3020.7
$ws.Range("M132").Value = -490.6999999999998
$ws.Range("H136").Value = 33925.41
$ws.Range("I136").Value = 39873.785
$ws.Range("K136").Value = 119621.355
$ws.Range("M136").Value = -117071.355
$ws.Range("H139").Value = 88015.8
$ws.Range("J139").Value = 116476.336
$ws.Range("L139").Value = 116476.336
$ws.Range("N139").Value = -126756.336

Write-Output "Applied 196 cell changes across 8 sheets"